$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 26 (weekly update), shifting
# the existing rows 26-36 down to 27-37.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with this week's reading.
$ws.Cells.Item(26, 1).Value = 2
$ws.Cells.Item(26, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44524
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100112032
$ws.Cells.Item(26, 7).Value = "Zapallo italiano"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 508
$ws.Cells.Item(26, 11).Value = 5000
$ws.Cells.Item(26, 12).Value = 6000
$ws.Cells.Item(26, 13).Value = 5508
$ws.Cells.Item(26, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 92
$ws.Cells.Item(26, 17).Value = 60
$ws.Cells.Item(26, 18).Value = "Hortaliza"
